$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 50.70817566666667
$ws.Range("H2").Value = 152.124527
$ws.Range("I2").Value = 0.5661129211027078
$ws.Range("J2").Value = 0.5661129211027077
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.99161333333333
$ws.Range("N2").Value = 92.97484
$ws.Range("O2").Value = 0.3599121977633812
$ws.Range("P2").Value = 0.3599121977633811
$ws.Range("Q2").Value = 1571.528173100076
$ws.Range("R2").Value = 14143.75355790068
$ws.Range("S2").Value = 0.2037509456163232
$ws.Range("T2").Value = 0.2037509456163231

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 50.70817566666667
$ws.Range("H3").Value = 152.124527
$ws.Range("I3").Value = 0.5661129211027078
$ws.Range("J3").Value = 0.5661129211027077
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 29.913269
$ws.Range("N3").Value = 89.739807
$ws.Range("O3").Value = 0.3473891556493311
$ws.Range("P3").Value = 0.3473891556493311
$ws.Range("Q3").Value = 1516.847299216254
$ws.Range("R3").Value = 13651.62569294629
$ws.Range("S3").Value = 0.1966614896640461
$ws.Range("T3").Value = 0.196661489664046

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 50.70817566666667
$ws.Range("H4").Value = 152.124527
$ws.Range("I4").Value = 0.5661129211027078
$ws.Range("J4").Value = 0.5661129211027077
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.150218
$ws.Range("N4").Value = 57.450654
$ws.Range("O4").Value = 0.2223955550134164
$ws.Range("P4").Value = 0.2223955550134163
$ws.Range("Q4").Value = 971.0726183989619
$ws.Range("R4").Value = 8739.653565590657
$ws.Range("S4").Value = 0.1259009972889031
$ws.Range("T4").Value = 0.1259009972889031

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 50.70817566666667
$ws.Range("H5").Value = 152.124527
$ws.Range("I5").Value = 0.5661129211027078
$ws.Range("J5").Value = 0.5661129211027077
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.053716000000001
$ws.Range("N5").Value = 18.161148
$ws.Range("O5").Value = 0.07030309157387134
$ws.Range("P5").Value = 0.07030309157387132
$ws.Range("Q5").Value = 306.9728943641107
$ws.Range("R5").Value = 2762.756049276996
$ws.Range("S5").Value = 0.03979948853343547
$ws.Range("T5").Value = 0.03979948853343545

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.08683666666667
$ws.Range("H6").Value = 51.26051
$ws.Range("I6").Value = 0.1907597520636141
$ws.Range("J6").Value = 0.1907597520636141
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.99161333333333
$ws.Range("N6").Value = 92.97484
$ws.Range("O6").Value = 0.3599121977633812
$ws.Range("P6").Value = 0.3599121977633811
$ws.Range("Q6").Value = 529.5486350631555
$ws.Range("R6").Value = 4765.937715568401
$ws.Range("S6").Value = 0.06865676161001304
$ws.Range("T6").Value = 0.06865676161001301

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.08683666666667
$ws.Range("H7").Value = 51.26051
$ws.Range("I7").Value = 0.1907597520636141
$ws.Range("J7").Value = 0.1907597520636141
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 29.913269
$ws.Range("N7").Value = 89.739807
$ws.Range("O7").Value = 0.3473891556493311
$ws.Range("P7").Value = 0.3473891556493311
$ws.Range("Q7").Value = 511.1231415690633
$ws.Range("R7").Value = 4600.108274121571
$ws.Range("S7").Value = 0.06626786920125464
$ws.Range("T7").Value = 0.06626786920125463

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.08683666666667
$ws.Range("H8").Value = 51.26051
$ws.Range("I8").Value = 0.1907597520636141
$ws.Range("J8").Value = 0.1907597520636141
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 19.150218
$ws.Range("N8").Value = 57.450654
$ws.Range("O8").Value = 0.2223955550134164
$ws.Range("P8").Value = 0.2223955550134163
$ws.Range("Q8").Value = 327.21664709706
$ws.Range("R8").Value = 2944.94982387354
$ws.Range("S8").Value = 0.04242412093440916
$ws.Range("T8").Value = 0.04242412093440914

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.08683666666667
$ws.Range("H9").Value = 51.26051
$ws.Range("I9").Value = 0.1907597520636141
$ws.Range("J9").Value = 0.1907597520636141
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.053716000000001
$ws.Range("N9").Value = 18.161148
$ws.Range("O9").Value = 0.07030309157387134
$ws.Range("P9").Value = 0.07030309157387132
$ws.Range("Q9").Value = 103.4388565183867
$ws.Range("R9").Value = 930.9497086654801
$ws.Range("S9").Value = 0.01341100031793725
$ws.Range("T9").Value = 0.01341100031793725

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 18.33915266666667
$ws.Range("H10").Value = 55.017458
$ws.Range("I10").Value = 0.2047407770084672
$ws.Range("J10").Value = 0.2047407770084672
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.99161333333333
$ws.Range("N10").Value = 92.97484
$ws.Range("O10").Value = 0.3599121977633812
$ws.Range("P10").Value = 0.3599121977633811
$ws.Range("Q10").Value = 568.3599283063022
$ws.Range("R10").Value = 5115.239354756721
$ws.Range("S10").Value = 0.07368870302489977
$ws.Range("T10").Value = 0.07368870302489974

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 18.33915266666667
$ws.Range("H11").Value = 55.017458
$ws.Range("I11").Value = 0.2047407770084672
$ws.Range("J11").Value = 0.2047407770084672
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 29.913269
$ws.Range("N11").Value = 89.739807
$ws.Range("O11").Value = 0.3473891556493311
$ws.Range("P11").Value = 0.3473891556493311
$ws.Range("Q11").Value = 548.5840069500673
$ws.Range("R11").Value = 4937.256062550607
$ws.Range("S11").Value = 0.0711247256519594
$ws.Range("T11").Value = 0.07112472565195939

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 18.33915266666667
$ws.Range("H12").Value = 55.017458
$ws.Range("I12").Value = 0.2047407770084672
$ws.Range("J12").Value = 0.2047407770084672
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 19.150218
$ws.Range("N12").Value = 57.450654
$ws.Range("O12").Value = 0.2223955550134164
$ws.Range("P12").Value = 0.2223955550134163
$ws.Range("Q12").Value = 351.198771501948
$ws.Range("R12").Value = 3160.788943517532
$ws.Range("S12").Value = 0.04553343873667618
$ws.Range("T12").Value = 0.04553343873667617

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 18.33915266666667
$ws.Range("H13").Value = 55.017458
$ws.Range("I13").Value = 0.2047407770084672
$ws.Range("J13").Value = 0.2047407770084672
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.053716000000001
$ws.Range("N13").Value = 18.161148
$ws.Range("O13").Value = 0.07030309157387134
$ws.Range("P13").Value = 0.07030309157387132
$ws.Range("Q13").Value = 111.0200219246427
$ws.Range("R13").Value = 999.1801973217841
$ws.Range("S13").Value = 0.01439390959493184
$ws.Range("T13").Value = 0.01439390959493183

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.438381
$ws.Range("H14").Value = 10.315143
$ws.Range("I14").Value = 0.03838654982521095
$ws.Range("J14").Value = 0.03838654982521095
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 30.99161333333333
$ws.Range("N14").Value = 92.97484
$ws.Range("O14").Value = 0.3599121977633812
$ws.Range("P14").Value = 0.3599121977633811
$ws.Range("Q14").Value = 106.56097444468
$ws.Range("R14").Value = 959.0487700021201
$ws.Range("S14").Value = 0.01381578751214521
$ws.Range("T14").Value = 0.01381578751214521

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.438381
$ws.Range("H15").Value = 10.315143
$ws.Range("I15").Value = 0.03838654982521095
$ws.Range("J15").Value = 0.03838654982521095
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 29.913269
$ws.Range("N15").Value = 89.739807
$ws.Range("O15").Value = 0.3473891556493311
$ws.Range("P15").Value = 0.3473891556493311
$ws.Range("Q15").Value = 102.853215777489
$ws.Range("R15").Value = 925.678941997401
$ws.Range("S15").Value = 0.01333507113207101
$ws.Range("T15").Value = 0.01333507113207101

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.438381
$ws.Range("H16").Value = 10.315143
$ws.Range("I16").Value = 0.03838654982521095
$ws.Range("J16").Value = 0.03838654982521095
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 19.150218
$ws.Range("N16").Value = 57.450654
$ws.Range("O16").Value = 0.2223955550134164
$ws.Range("P16").Value = 0.2223955550134163
$ws.Range("Q16").Value = 65.84574571705799
$ws.Range("R16").Value = 592.611711453522
$ws.Range("S16").Value = 0.008536998053427952
$ws.Range("T16").Value = 0.00853699805342795

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.438381
$ws.Range("H17").Value = 10.315143
$ws.Range("I17").Value = 0.03838654982521095
$ws.Range("J17").Value = 0.03838654982521095
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.053716000000001
$ws.Range("N17").Value = 18.161148
$ws.Range("O17").Value = 0.07030309157387134
$ws.Range("P17").Value = 0.07030309157387132
$ws.Range("Q17").Value = 20.814982073796
$ws.Range("R17").Value = 187.334838664164
$ws.Range("S17").Value = 0.00269869312756678
$ws.Range("T17").Value = 0.00269869312756678

